# Multiple changes:
# 1) Clean up the Saleccia "Commentaires" note (F3): replace the old
#    hyphen-bulleted, line-broken text with a single sentence-per-period
#    version (no more "&#10;" line break, no leading "-" bullets).
# 2) The sheet that is active/selected moves from "Lotu" to "Saleccia",
#    with the selection landing on F3 (the cell that was just cleaned up).

$wb = $excel.ActiveWorkbook

$saleccia = $wb.Worksheets.Item("Saleccia")

# Replace the comment text in F3 - drop the leading "-" bullets and the
# embedded line break, turning it into two clean sentences.
$saleccia.Range("F3").Value = "Bateaux passant dans la zone de baignade. Bateaux dans le chenal non motorisé."

# Make Saleccia the active sheet/tab and put the selection on F3.
$saleccia.Activate() | Out-Null
$saleccia.Range("F3").Select() | Out-Null
